# This script rewrites the full data table on Sheet1 (rows 2-34, columns A-F)
# to match the updated dataset described in the commit message:
#   "figuras atualizadas manualmente antes de criar definir a atualizacao automatica"
# Each region (Brasil, Nordeste, Sergipe) gains a new final year (2025) of data,
# which pushes the blocks for the later regions further down the sheet, and the
# sheet dimension grows from A1:F31 to A1:F34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One entry per data row: Row, Regiao, Ano, Variavel, Valor, Posicao (or $null), Faltam dados (bool)
$data = @(
    @{ Row = 2; A = "Brasil"; B = "01/01/2015"; C = "Estupro"; D = 23.51597416558056; E = $null; F = $true },
    @{ Row = 3; A = "Brasil"; B = "01/01/2016"; C = "Estupro"; D = 29.68688409921749; E = $null; F = $true },
    @{ Row = 4; A = "Brasil"; B = "01/01/2017"; C = "Estupro"; D = 31.95310127371981; E = $null; F = $true },
    @{ Row = 5; A = "Brasil"; B = "01/01/2018"; C = "Estupro"; D = 40.39619814149729; E = $null; F = $false },
    @{ Row = 6; A = "Brasil"; B = "01/01/2019"; C = "Estupro"; D = 43.73582517596133; E = $null; F = $false },
    @{ Row = 7; A = "Brasil"; B = "01/01/2020"; C = "Estupro"; D = 31.08082262272347; E = $null; F = $true },
    @{ Row = 8; A = "Brasil"; B = "01/01/2021"; C = "Estupro"; D = 40.37773101085241; E = $null; F = $false },
    @{ Row = 9; A = "Brasil"; B = "01/01/2022"; C = "Estupro"; D = 35.68146378041867; E = $null; F = $false },
    @{ Row = 10; A = "Brasil"; B = "01/01/2023"; C = "Estupro"; D = 41.76481125506481; E = $null; F = $false },
    @{ Row = 11; A = "Brasil"; B = "01/01/2024"; C = "Estupro"; D = 38.22179512271354; E = $null; F = $false },
    @{ Row = 12; A = "Brasil"; B = "01/01/2025"; C = "Estupro"; D = 13.22637324315908; E = $null; F = $false },
    @{ Row = 13; A = "Nordeste"; B = "01/01/2015"; C = "Estupro"; D = 12.70230727830731; E = $null; F = $true },
    @{ Row = 14; A = "Nordeste"; B = "01/01/2016"; C = "Estupro"; D = 14.63737214627904; E = $null; F = $true },
    @{ Row = 15; A = "Nordeste"; B = "01/01/2017"; C = "Estupro"; D = 14.84706792354145; E = $null; F = $true },
    @{ Row = 16; A = "Nordeste"; B = "01/01/2018"; C = "Estupro"; D = 15.02477095090394; E = $null; F = $false },
    @{ Row = 17; A = "Nordeste"; B = "01/01/2019"; C = "Estupro"; D = 17.86668182339003; E = $null; F = $false },
    @{ Row = 18; A = "Nordeste"; B = "01/01/2020"; C = "Estupro"; D = 16.36815773591297; E = $null; F = $true },
    @{ Row = 19; A = "Nordeste"; B = "01/01/2021"; C = "Estupro"; D = 18.080423162261; E = $null; F = $false },
    @{ Row = 20; A = "Nordeste"; B = "01/01/2022"; C = "Estupro"; D = 17.06753006920804; E = $null; F = $false },
    @{ Row = 21; A = "Nordeste"; B = "01/01/2023"; C = "Estupro"; D = 18.44416896158795; E = $null; F = $false },
    @{ Row = 22; A = "Nordeste"; B = "01/01/2024"; C = "Estupro"; D = 20.42957623803085; E = $null; F = $false },
    @{ Row = 23; A = "Nordeste"; B = "01/01/2025"; C = "Estupro"; D = 6.773419516901637; E = $null; F = $false },
    @{ Row = 24; A = "Sergipe"; B = "01/01/2015"; C = "Estupro"; D = 0.0; E = 25.0; F = $true },
    @{ Row = 25; A = "Sergipe"; B = "01/01/2016"; C = "Estupro"; D = 0.0; E = 25.5; F = $true },
    @{ Row = 26; A = "Sergipe"; B = "01/01/2017"; C = "Estupro"; D = 0.0; E = 25.5; F = $true },
    @{ Row = 27; A = "Sergipe"; B = "01/01/2018"; C = "Estupro"; D = 1.529031205827987; E = 27.0; F = $false },
    @{ Row = 28; A = "Sergipe"; B = "01/01/2019"; C = "Estupro"; D = 12.70379708923197; E = 23.0; F = $false },
    @{ Row = 29; A = "Sergipe"; B = "01/01/2020"; C = "Estupro"; D = 11.91875244834513; E = 19.0; F = $true },
    @{ Row = 30; A = "Sergipe"; B = "01/01/2021"; C = "Estupro"; D = 11.72858727514357; E = 22.0; F = $false },
    @{ Row = 31; A = "Sergipe"; B = "01/01/2022"; C = "Estupro"; D = 13.67300237025999; E = 21.0; F = $false },
    @{ Row = 32; A = "Sergipe"; B = "01/01/2023"; C = "Estupro"; D = 15.99320654049082; E = 19.0; F = $false },
    @{ Row = 33; A = "Sergipe"; B = "01/01/2024"; C = "Estupro"; D = 13.9306383050681; E = 23.0; F = $false },
    @{ Row = 34; A = "Sergipe"; B = "01/01/2025"; C = "Estupro"; D = 9.267759143722971; E = 11.0; F = $false }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A             # Regiao
    $ws.Cells.Item($r, 2).Value = "'" + $item.B       # Ano (apostrophe forces text, avoids date coercion)
    $ws.Cells.Item($r, 3).Value = $item.C             # Variavel
    $ws.Cells.Item($r, 4).Value = $item.D             # Valor
    if ($null -eq $item.E) {
        $ws.Cells.Item($r, 5).Value = ""              # Posicao (blank)
    } else {
        $ws.Cells.Item($r, 5).Value = $item.E         # Posicao relativamente as demais UF
    }
    $ws.Cells.Item($r, 6).Value = $item.F             # Faltam dados para todos os Estados
}
